$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$layout4 = $master.CustomLayouts.Item(4)
$content1 = $layout4.Shapes.Item("Content Placeholder 2")
$content1.Width = 7992000 / 12700
Write-Host "Width prop: $($content1.Width)"
